$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 278199
$ws.Range("E10").Value = 1752276968
$ws.Range("C19").Value = 108920
$ws.Range("E19").Value = 344665282
$ws.Range("C115").Value = 17556
$ws.Range("E115").Value = 38619107
$ws.Range("C134").Value = 5678
$ws.Range("E134").Value = 17162226
$ws.Range("C152").Value = 126047
$ws.Range("E152").Value = 715970384
$ws.Range("C168").Value = 285023
$ws.Range("E168").Value = 1210861599
$ws.Range("C169").Value = 562613
$ws.Range("E169").Value = 1285058739
$ws.Range("C170").Value = 367415
$ws.Range("E170").Value = 2846304604
$ws.Range("C171").Value = 115169
$ws.Range("E171").Value = 447140114
$ws.Range("C173").Value = 54392
$ws.Range("E173").Value = 151908108
$ws.Range("C174").Value = 357255
$ws.Range("E174").Value = 1018552234
$ws.Range("C175").Value = 125559
$ws.Range("E175").Value = 813138243
$ws.Range("C177").Value = 96761
$ws.Range("E177").Value = 174754897
$ws.Range("C179").Value = 235723
$ws.Range("E179").Value = 812710859
$ws.Range("C267").Value = 84977
$ws.Range("E267").Value = 156521773
$ws.Range("C313").Value = 220650
$ws.Range("E313").Value = 1371052430
$ws.Range("C317").Value = 103584
$ws.Range("E317").Value = 303300964
